# Rename the column header cells so the "_old"/"_new" suffixes become the
# respective format-version identifiers ("_FV2404" / "_FV2410"), then turn
# the header + data range into a real Excel Table, and finally freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename headers: "<name>_old" -> "<name>_FV2404" (A1:J1) ---------
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"

# column K1 ("diff") is unchanged

# --- 2. Rename headers: "<name>_new" -> "<name>_FV2410" (L1:U1) ---------
$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- 3. Turn A1:U54 into an Excel Table (adds xl/tables/table1.xml) -----
$tableRange = $ws.Range("A1:U54")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 4. Freeze the header row (row 1) ------------------------------------
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
